# Daily attendance processing - 2025-12-13 01:26:51
# Rotate the comma-separated "Recorded By" list in column G so that the
# last listed recorder is moved to the front of the list, unless the
# last entry is already the exact, case-sensitive literal "System".

function Test-ExactEquals($x, $y) {
    if ($x.Length -ne $y.Length) { return $false }
    $cx = $x.ToCharArray()
    $cy = $y.ToCharArray()
    for ($i = 0; $i -lt $cx.Length; $i++) {
        if ([int]$cx[$i] -ne [int]$cy[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"
    if ($parts.Count -gt 1 -and -not (Test-ExactEquals $parts[-1] "System")) {
        $newParts = @($parts[-1]) + $parts[0..($parts.Count - 2)]
        $newValue = $newParts -join ", "
        $cell.Value = $newValue
    }
}
